$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("I2").Value = 3651
$ws.Range("I3").Value = 3773
$ws.Range("G4").Value = 1436
$ws.Range("I4").Value = 883
$ws.Range("I6").Value = 4233
$ws.Range("G7").Value = 24658
$ws.Range("I7").Value = 12888

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("I7").Value = 417
$ws.Range("I8").Value = 783
$ws.Range("G10").Value = 189
$ws.Range("I11").Value = 203
$ws.Range("I14").Value = 65
$ws.Range("I15").Value = 152
$ws.Range("I19").Value = 343
$ws.Range("I20").Value = 317
$ws.Range("I23").Value = 124
$ws.Range("I27").Value = 120
$ws.Range("I29").Value = 837
$ws.Range("I31").Value = 121
$ws.Range("I33").Value = 580
$ws.Range("I36").Value = 178
$ws.Range("I37").Value = 413
$ws.Range("I42").Value = 451
$ws.Range("I47").Value = 90
$ws.Range("I53").Value = 142
$ws.Range("I54").Value = 289
$ws.Range("I56").Value = 14
$ws.Range("I58").Value = 9
$ws.Range("I59").Value = 27
$ws.Range("I63").Value = 49
$ws.Range("I64").Value = 115
$ws.Range("I67").Value = 498
$ws.Range("I72").Value = 48
$ws.Range("I75").Value = 46
$ws.Range("I76").Value = 197
$ws.Range("I79").Value = 342
$ws.Range("I83").Value = 259
$ws.Range("I85").Value = 586
$ws.Range("I88").Value = 117
$ws.Range("I89").Value = 146
$ws.Range("I92").Value = 39
$ws.Range("I93").Value = 73
$ws.Range("I95").Value = 203
$ws.Range("I96").Value = 144
$ws.Range("I99").Value = 240
$ws.Range("G101").Value = 24658
$ws.Range("I101").Value = 12888

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("I3").Value = 235
$ws.Range("I4").Value = 35
$ws.Range("I7").Value = 586

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("I2").Value = 92
$ws.Range("I7").Value = 203

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("I2").Value = 244
$ws.Range("I6").Value = 253
$ws.Range("I7").Value = 783

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range("I2").Value = 29
$ws.Range("I3").Value = 37
$ws.Range("I7").Value = 142

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("I2").Value = 144
$ws.Range("I7").Value = 417

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range("I2").Value = 37
$ws.Range("I3").Value = 32
$ws.Range("I7").Value = 146

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range("I3").Value = 43
$ws.Range("I7").Value = 144

$ws = $wb.Worksheets.Item('Bridgeport')
$ws.Range("I2").Value = 19
$ws.Range("I7").Value = 65

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("I3").Value = 126
$ws.Range("I7").Value = 413

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range("I6").Value = 67
$ws.Range("I7").Value = 240

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("I2").Value = 120
$ws.Range("I3").Value = 174
$ws.Range("I6").Value = 167
$ws.Range("I7").Value = 498

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range("I2").Value = 40
$ws.Range("I7").Value = 121

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("I2").Value = 90
$ws.Range("I3").Value = 101
$ws.Range("I7").Value = 259

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("I3").Value = 76
$ws.Range("I7").Value = 203

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("I3").Value = 210
$ws.Range("I6").Value = 185
$ws.Range("I7").Value = 580

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("I2").Value = 65
$ws.Range("I7").Value = 289

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("I3").Value = 287
$ws.Range("I6").Value = 227
$ws.Range("I7").Value = 837

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("I2").Value = 132
$ws.Range("I7").Value = 343

$ws = $wb.Worksheets.Item('River North')
$ws.Range("I2").Value = 42
$ws.Range("I7").Value = 197

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("I3").Value = 153
$ws.Range("I7").Value = 451

$ws = $wb.Worksheets.Item('Avondale')
$ws.Range("G4").Value = 8
$ws.Range("G7").Value = 189

$ws = $wb.Worksheets.Item('Douglas')
$ws.Range("I2").Value = 34
$ws.Range("I7").Value = 124

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("I4").Value = 19
$ws.Range("I7").Value = 342

$ws = $wb.Worksheets.Item('Near South Side')
$ws.Range("I2").Value = 30
$ws.Range("I7").Value = 115

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("I4").Value = 21
$ws.Range("I7").Value = 317

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range("I2").Value = 57
$ws.Range("I3").Value = 54
$ws.Range("I7").Value = 178

$ws = $wb.Worksheets.Item('West Lawn')
$ws.Range("I2").Value = 23
$ws.Range("I7").Value = 73

$ws = $wb.Worksheets.Item('Kenwood')
$ws.Range("I3").Value = 30
$ws.Range("I7").Value = 90

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range("I3").Value = 35
$ws.Range("I4").Value = 11
$ws.Range("I7").Value = 152

$ws = $wb.Worksheets.Item('Montclare')
$ws.Range("I2").Value = 12
$ws.Range("I7").Value = 27

$ws = $wb.Worksheets.Item('West Elsdon')
$ws.Range("I3").Value = 8
$ws.Range("I7").Value = 39

$ws = $wb.Worksheets.Item('United Center')
$ws.Range("I3").Value = 41
$ws.Range("I7").Value = 117

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range("I4").Value = 17
$ws.Range("I6").Value = 49
$ws.Range("I7").Value = 120

$ws = $wb.Worksheets.Item('Pullman')
$ws.Range("I3").Value = 16
$ws.Range("I4").Value = 2
$ws.Range("I6").Value = 11
$ws.Range("I7").Value = 46

$ws = $wb.Worksheets.Item('Old Town')
$ws.Range("I2").Value = 8
$ws.Range("I7").Value = 48

$ws = $wb.Worksheets.Item('Magnificent Mile')
$ws.Range("I2").Value = 4
$ws.Range("I7").Value = 14

$ws = $wb.Worksheets.Item('Archer Heights')
$ws.Range("I4").Value = 1

$ws = $wb.Worksheets.Item('Millenium Park')
$ws.Range("I7").Value = 9
